$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.175.63"
$ws.Range("E2").Value = "  +5.34%  "

$ws.Range("D3").Value = "1.783.25"
$ws.Range("E3").Value = "  +2.96%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.51"
$ws.Range("E5").Value = "  +0.78%  "

$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4912"
$ws.Range("E7").Value = "  -0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2678"
$ws.Range("E8").Value = "  +1.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06271"
$ws.Range("E9").Value = "  +0.74%  "

$ws.Range("D10").Value = "1.773.24"
$ws.Range("E10").Value = "  +2.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.48"
$ws.Range("E11").Value = "  +3.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07036"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6284"
$ws.Range("E13").Value = "  +2.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.656"
$ws.Range("E14").Value = "  +3.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "80.04"
$ws.Range("E15").Value = "  +3.50%  "

$ws.Range("D16").Value = "28.150.15"
$ws.Range("E16").Value = "  +6.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007240"
$ws.Range("E19").Value = "  -0.11%  "

$ws.Range("E20").Value = "  +5.21%  "

$ws.Range("D21").Value = "2.008.83"
$ws.Range("E21").Value = "  +3.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.555"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.741"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.270"
$ws.Range("E24").Value = "  +3.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.18"
$ws.Range("E25").Value = "  +2.10%  "

$ws.Range("E26").Value = "  +2.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.859"
$ws.Range("E27").Value = "  +4.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.71"
$ws.Range("E28").Value = "  +2.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.385"
$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.193"
$ws.Range("E30").Value = "  +6.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08259"
$ws.Range("E31").Value = "  +3.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.767"
$ws.Range("E32").Value = "  +2.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04894"
$ws.Range("E33").Value = "  +9.15%  "

$ws.Range("E34").Value = "  +7.13%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.617"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6517"
$ws.Range("E36").Value = "  +4.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9505"
$ws.Range("E37").Value = "  +1.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.605"
$ws.Range("E38").Value = "  +7.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.043"
$ws.Range("E39").Value = "  -0.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.914"
$ws.Range("E40").Value = "  +5.27%  "

$ws.Range("E41").Value = "  +2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9994"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.88"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3983"
$ws.Range("E44").Value = "  +3.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.181"
$ws.Range("E45").Value = "  +3.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1219"
$ws.Range("E46").Value = "  +5.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05443"
$ws.Range("E47").Value = "  +1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.987"
$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("E49").Value = "  +4.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.75"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.92"
$ws.Range("E51").Value = "  +2.11%  "
